$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# A new weekly price observation was inserted as a data row right before the
# existing row 307, pushing row 307 and every row after it down by one
# (old row 307 -> new row 308, ..., old row 413 -> new row 414).
$ws.Rows.Item(307).Insert()

# Populate the newly inserted row 307 with the new observation's data. The
# descriptive columns mirror the surrounding rows for this market/category.
$ws.Cells.Item(307, 1).Value = 4
$ws.Cells.Item(307, 2).Value = 'Feria Lagunitas de Puerto Montt'
$ws.Cells.Item(307, 3).Value = 'Los Lagos'
$ws.Cells.Item(307, 4).Value = 44988
$ws.Cells.Item(307, 5).Value = 10
$ws.Cells.Item(307, 6).Value = 100112037
$ws.Cells.Item(307, 7).Value = 'Cebollín'
$ws.Cells.Item(307, 8).Value = 'Sin especificar'
$ws.Cells.Item(307, 9).Value = 'Primera'
$ws.Cells.Item(307, 10).Value = 140
$ws.Cells.Item(307, 11).Value = 7500
$ws.Cells.Item(307, 12).Value = 7500
$ws.Cells.Item(307, 13).Value = 7500
$ws.Cells.Item(307, 14).Value = '$/paquete 36 unidades'
$ws.Cells.Item(307, 15).Value = 'Región Metropolitana'
$ws.Cells.Item(307, 16).Value = 208
$ws.Cells.Item(307, 17).Value = 36
$ws.Cells.Item(307, 18).Value = 'Hortaliza'
